# Append two new test-case rows (RS_10544 / Source Object Update) to the
# TestData sheet, mirroring the existing row-pair pattern (header-ish row
# followed by a "values" row) already used throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18/19 is an existing pair that uses the exact same row- and cell-
# level style pairing as the new rows need (A/D/E bold-ish label style, B
# "search" style, C "process name" style for row 26; plain value style
# for row 27), so copy their formatting down to the new rows before
# writing values.
$ws.Range("A18:E18").Copy($ws.Range("A26"))
$ws.Range("A19:E19").Copy($ws.Range("A27"))
$ws.Rows(26).RowHeight = $ws.Rows(18).RowHeight
$ws.Rows(27).RowHeight = $ws.Rows(19).RowHeight

# Row 26 - "SCN_SourceObjectUpdate_RS_10544" scenario definition
$ws.Range("A26").Value = "SCN_SourceObjectUpdate_RS_10544"
$ws.Range("B26").Value = "ExploreSearch"
$ws.Range("C26").Value = "ProcessName"
$ws.Range("D26").Value = "ExploreChildSearch"
$ws.Range("E26").Value = "EditProcessName"

# Row 27 - matching data/value row
$ws.Range("A27").Value = "Data_SCN_SourceObjectUpdate_RS_10544"
$ws.Range("B27").Value = "AUTOMATION SEARCH"
$ws.Range("C27").Value = "Manual_RS_10544_SOU"
$ws.Range("D27").Value = "Work Orders"
$ws.Range("E27").Value = "Manual_Checklist_Edit_WO_PROCESS"

# Reflect the author's final cursor/scroll position (B32, scrolled so row
# 11 is the top visible row).
[void]$ws.Range("B32").Select()
$excel.ActiveWindow.ScrollRow = 11
